$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order: B C D E F I J K L M N (G stays 1, H stays empty)
$cols = 2,3,4,5,6,9,10,11,12,13,14

$data = @{
    2 = @(1.02, 1.018067105082657, 1.028050429999862, 1.018834723825937, 1.028569057437104, 1.029841315527923, 1.023277645916982, 1.030868375034837, 1.021679761590334, 1.031385493480437, 1.024730817545772)
    3 = @(1.02, 1.022104602649436, 1.028550325515335, 1.022395607626298, 1.032018474363051, 1.030032628188528, 1.026935162657732, 1.031176956395566, 1.025038977273434, 1.034635768186816, 1.028393528379816)
    4 = @(1.02, 1.024687327660341, 1.028871716451234, 1.024672132994881, 1.034223009948832, 1.030152752987494, 1.029273053295223, 1.031373842521862, 1.027185086339882, 1.036711476794205, 1.030734739090171)
    5 = @(1.02, 1.02576621780334, 1.029006342577223, 1.02562279538125, 1.035143425828342, 1.030202387868978, 1.030249243701315, 1.031455956792303, 1.028080932488834, 1.037577734274358, 1.031712315798611)
    6 = @(1.02, 1.025946971380876, 1.029028918646499, 1.025782047382171, 1.035297599944332, 1.030210671376977, 1.030412766378771, 1.031469705929866, 1.028230981019669, 1.037722814782982, 1.031876070697025)
    7 = @(1.02, 1.024701770582887, 1.028873517230918, 1.024684860609804, 1.034235333356784, 1.030153419596169, 1.029286123066528, 1.031374942303238, 1.027197081465032, 1.036723076542353, 1.030747827422051)
    8 = @(1.02, 1.01943791994067, 1.028219807873862, 1.020043991594583, 1.029740617097775, 1.029906737470805, 1.024519817150466, 1.030973244983417, 1.022820854235401, 1.032489741293921, 1.025974752804885)
    9 = @(1.02, 1.009922495126758, 1.027051520997268, 1.01164452576869, 1.021600426788542, 1.029443385824711, 1.015890017250645, 1.030243568875564, 1.014888656369322, 1.024810704611762, 1.017332697599354)
    10 = @(1.02, 1.003400783643963, 1.02626095705391, 1.005880944870363, 1.016011866817545, 1.029114364405802, 1.009966054863218, 1.029741716813135, 1.009437723447311, 1.019530434441688, 1.011400322505849)
    11 = @(1.02, 1.000530614495827, 1.025915690400758, 1.003342866746604, 1.01355032597436, 1.028966922186335, 1.007356775589569, 1.029520585223876, 1.00703541340947, 1.01720269938594, 1.008787337756343)
    12 = @(1.02, 0.9994572061126786, 1.025786984999676, 1.002393425777833, 1.012629448606412, 1.028911390262416, 1.006380611505209, 1.029437856740018, 1.006136471869503, 1.016331580505752, 1.007809787407016)
    13 = @(1.02, 0.9996877920790229, 1.025814613723718, 1.002597391849281, 1.012827281259715, 1.02892333698123, 1.006590322460577, 1.029455629253989, 1.006329602338462, 1.016518737049293, 1.008019796175995)
    14 = @(1.02, 1.000442037251159, 1.025905061005926, 1.003264523918095, 1.013474341356503, 1.028962347620145, 1.00727622951968, 1.029513758997395, 1.006961243406443, 1.017130826597673, 1.008706677301796)
    15 = @(1.02, 1.000905774392833, 1.025960727404263, 1.003674670420382, 1.013872139298286, 1.028986281394284, 1.007697906986828, 1.02954949594415, 1.007349532216102, 1.017507085893554, 1.009128953599315)
    16 = @(1.02, 1.003590263504915, 1.026283807890807, 1.00604846876641, 1.016174328820483, 1.029124043347452, 1.010138266145564, 1.029756310667316, 1.009596245961826, 1.019684023646646, 1.011572778347969)
    17 = @(1.02, 1.00526155166576, 1.026485668210176, 1.007525918485579, 1.017607076223123, 1.029209113692686, 1.011656989460346, 1.029885004971285, 1.010994090467136, 1.021038294873102, 1.013093658424058)
    18 = @(1.02, 1.006231950592588, 1.026603126623141, 1.008383620007429, 1.018438774710509, 1.02925825515632, 1.012538596704218, 1.029959702085257, 1.01180539644898, 1.0218242491694, 1.013976517651338)
    19 = @(1.02, 1.006562089759534, 1.026643129302407, 1.008675393024243, 1.01872169279206, 1.029274930475086, 1.012838492739044, 1.029985109892511, 1.012081355719511, 1.022091574398898, 1.014276839572933)
    20 = @(1.02, 1.005082699525005, 1.026464039918263, 1.007367825304419, 1.017453771676311, 1.029200036084337, 1.011494485524475, 1.029871235473854, 1.01084453428702, 1.020893406942305, 1.012930923713959)
    21 = @(1.02, 1.000220135295668, 1.025878439304846, 1.003068257312265, 1.013283981563855, 1.02895088123559, 1.00707444186269, 1.029496657662444, 1.006775426020084, 1.01695076310341, 1.008504603083189)
    22 = @(1.02, 0.9971204252406206, 1.025507591001268, 1.000326099606439, 1.010624208438364, 1.028789790168186, 1.004254933905137, 1.029257721768175, 1.004178573002791, 1.014434141199648, 1.005681091100933)
    23 = @(1.02, 0.9987677845559791, 1.025704441945322, 1.001783560689006, 1.012037913797515, 1.028875614492979, 1.005753555945905, 1.02938471615154, 1.005558963207922, 1.015771924194063, 1.007181841356894)
    24 = @(1.02, 1.005163528789978, 1.026473813686994, 1.007439273363141, 1.017523055739231, 1.029204139342646, 1.011567927133131, 1.029877458453562, 1.010912124729309, 1.020958887807798, 1.013004469618123)
    25 = @(1.02, 1.012412569937515, 1.027355560703707, 1.01384375002125, 1.023732316058226, 1.029566655092902, 1.018149947662938, 1.030434867109709, 1.016966918933758, 1.016518737049293, 1.019595837372051)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $vals[$i]
    }
}
